$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.644.14"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "2.673.77"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "2.673.17"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.360"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "3.164.56"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "67.593.86"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").Value = "2.672.31"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  -3.64%  "
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("E24").Value = "  -4.12%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.20%  "
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "557.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("E34").Value = "  -0.85%  "
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.373"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.50%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("E47").Value = "  -5.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.591"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.83%  "
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.76%  "
